# Auto-generated Excel COM-interop script
# Updates LeveHub market-price snapshot values (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the ALC, ARM, BSM, CRP, CUL, GSM and WVR sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 25057306
$ws.Range("I6").Value = 33335836
$ws.Range("K6").Value = 100007508
$ws.Range("M6").Value = -100007396
$ws.Range("H40").Value = 1187.8823
$ws.Range("I40").Value = 811.64703
$ws.Range("J40").Value = 1376
$ws.Range("K40").Value = 811.64703
$ws.Range("L40").Value = 1376
$ws.Range("M40").Value = -636.64703
$ws.Range("N40").Value = -1726
$ws.Range("H100").Value = 748.1667
$ws.Range("I100").Value = 822.25
$ws.Range("J100").Value = 600
$ws.Range("K100").Value = 822.25
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = -281.25
$ws.Range("N100").Value = -1682
$ws.Range("H112").Value = 26317194
$ws.Range("I112").Value = 333333730
$ws.Range("J112").Value = 1491.7715
$ws.Range("K112").Value = 1000001190
$ws.Range("L112").Value = 4475.3145
$ws.Range("M112").Value = -1000000082
$ws.Range("N112").Value = -6691.3145
$ws.Range("H118").Value = 1355.5555
$ws.Range("I118").Value = 966.6667
$ws.Range("J118").Value = 1550
$ws.Range("K118").Value = 2900.0001
$ws.Range("L118").Value = 4650
$ws.Range("N118").Value = -7964
$ws.Range("M118").Value = -1243.0001
$ws.Range("H129").Value = 897.9792
$ws.Range("J129").Value = 992.1053000000001
$ws.Range("L129").Value = 2976.3159
$ws.Range("N129").Value = -12976.3159
$ws.Range("H132").Value = 1325855.5
$ws.Range("I132").Value = 1513.0588
$ws.Range("J132").Value = 16335070
$ws.Range("K132").Value = 4539.1764
$ws.Range("L132").Value = 49005210
$ws.Range("M132").Value = -2009.1764
$ws.Range("N132").Value = -49010270
$ws.Range("H133").Value = 29915.385
$ws.Range("J133").Value = 29915.385
$ws.Range("L133").Value = 29915.385
$ws.Range("N133").Value = -40035.38499999999
$ws.Range("H136").Value = 49780
$ws.Range("J136").Value = 49780
$ws.Range("L136").Value = 49780
$ws.Range("M136").Value = -59980
$ws.Range("H137").Value = 2779405.8
$ws.Range("I137").Value = 3126464.5
$ws.Range("J137").Value = 2937.5
$ws.Range("K137").Value = 9379393.5
$ws.Range("L137").Value = 8812.5
$ws.Range("M137").Value = -9376843.5
$ws.Range("N137").Value = -13912.5
$ws.Range("H138").Value = 2490011.8
$ws.Range("I138").Value = 1345.5186
$ws.Range("J138").Value = 4169861.5
$ws.Range("K138").Value = 4036.5558
$ws.Range("L138").Value = 12509584.5
$ws.Range("N138").Value = -12519864.5
$ws.Range("M138").Value = 1103.4442
$ws.Range("H140").Value = 63004.4
$ws.Range("J140").Value = 63004.4
$ws.Range("L140").Value = 63004.4
$ws.Range("N140").Value = -73364.39999999999
$ws.Range("H141").Value = 732.3333
$ws.Range("I141").Value = 732.3333
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2196.9999
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2983.0001
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7001799
$ws.Range("I74").Value = 9655164
$ws.Range("J74").Value = 103047.8
$ws.Range("K74").Value = 9655164
$ws.Range("L74").Value = 103047.8
$ws.Range("M74").Value = -9654290
$ws.Range("N74").Value = -104795.8
$ws.Range("H77").Value = 7001799
$ws.Range("I77").Value = 9655164
$ws.Range("J77").Value = 103047.8
$ws.Range("K77").Value = 48275820
$ws.Range("L77").Value = 515239
$ws.Range("M77").Value = -48271452
$ws.Range("N77").Value = -523975
$ws.Range("H110").Value = 1251903.8
$ws.Range("I110").Value = 1429461.4
$ws.Range("J110").Value = 9000
$ws.Range("K110").Value = 1429461.4
$ws.Range("L110").Value = 9000
$ws.Range("M110").Value = -1427416.4
$ws.Range("N110").Value = -13090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1092.7142
$ws.Range("I94").Value = 1129.8
$ws.Range("K94").Value = 1129.8
$ws.Range("M94").Value = -678.8
$ws.Range("H99").Value = 1042
$ws.Range("I99").Value = 1071.25
$ws.Range("J99").Value = 925
$ws.Range("K99").Value = 1071.25
$ws.Range("L99").Value = 925
$ws.Range("M99").Value = 426.75
$ws.Range("N99").Value = -3921
$ws.Range("H105").Value = 41668584
$ws.Range("I105").Value = 50001880
$ws.Range("J105").Value = 2111
$ws.Range("K105").Value = 50001880
$ws.Range("L105").Value = 2111
$ws.Range("M105").Value = -50000133
$ws.Range("N105").Value = -5605
$ws.Range("H107").Value = 1561.6428
$ws.Range("I107").Value = 1687.4546
$ws.Range("K107").Value = 1687.4546
$ws.Range("M107").Value = 232.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 83334410
$ws.Range("I16").Value = 1190
$ws.Range("J16").Value = 111112140
$ws.Range("K16").Value = 1190
$ws.Range("L16").Value = 111112140
$ws.Range("M16").Value = -903
$ws.Range("N16").Value = -111112714
$ws.Range("H113").Value = 83334410
$ws.Range("I113").Value = 1190
$ws.Range("J113").Value = 111112140
$ws.Range("K113").Value = 1190
$ws.Range("L113").Value = 111112140
$ws.Range("M113").Value = 980
$ws.Range("N113").Value = -111116480
$ws.Range("H134").Value = 32966.113
$ws.Range("I134").Value = 2377.7778
$ws.Range("J134").Value = 136201.75
$ws.Range("K134").Value = 7133.3334
$ws.Range("L134").Value = 408605.25
$ws.Range("M134").Value = -4598.3334
$ws.Range("N134").Value = -413675.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 505.4091
$ws.Range("I7").Value = 135.44444
$ws.Range("J7").Value = 761.53845
$ws.Range("K7").Value = 406.33332
$ws.Range("L7").Value = 2284.61535
$ws.Range("M7").Value = -294.33332
$ws.Range("N7").Value = -2508.61535
$ws.Range("H80").Value = 3983.3333
$ws.Range("J80").Value = 3981.818
$ws.Range("L80").Value = 11945.454
$ws.Range("N80").Value = -13817.454
$ws.Range("H83").Value = 3983.3333
$ws.Range("J83").Value = 3981.818
$ws.Range("L83").Value = 35836.362
$ws.Range("N83").Value = -45196.362
$ws.Range("H92").Value = 827.0909
$ws.Range("I92").Value = 827.0909
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2481.2727
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1233.2727
$ws.Range("N92").ClearContents()
$ws.Range("H131").Value = 997.4
$ws.Range("I131").Value = 575.4286
$ws.Range("J131").Value = 1037.863
$ws.Range("K131").Value = 1726.2858
$ws.Range("L131").Value = 3113.589
$ws.Range("M131").Value = 3313.7142
$ws.Range("N131").Value = -13193.589

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1480
$ws.Range("I113").Value = 1300
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 1300
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = 870
$ws.Range("N113").Value = -6540

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 279.73685
$ws.Range("J107").Value = 294.27274
$ws.Range("L107").Value = 882.81822
$ws.Range("N107").Value = -4722.81822
